$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1451232.1
$ws.Range("J17").Value = 1494532.1
$ws.Range("L17").Value = 4483596.300000001
$ws.Range("N17").Value = -4483932.300000001

$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652
$ws.Range("M32").ClearContents()

$ws.Range("H132").Value = 9474
$ws.Range("I132").Value = 9474
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 28422
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -25892
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 48975
$ws.Range("J134").Value = 48975
$ws.Range("L134").Value = 48975
$ws.Range("N134").Value = -59115

$ws.Range("H135").Value = 35724650
$ws.Range("I135").Value = 1008.4
$ws.Range("K135").Value = 9075.6
$ws.Range("M135").Value = -6540.6

$ws.Range("H137").Value = 2426.3333
$ws.Range("I137").Value = 2373.2273
$ws.Range("J137").Value = 2660
$ws.Range("K137").Value = 7119.6819
$ws.Range("L137").Value = 7980
$ws.Range("M137").Value = -4569.6819
$ws.Range("N137").Value = -13080

$ws.Range("H138").Value = 1570.5927
$ws.Range("J138").Value = 2279.8696
$ws.Range("L138").Value = 6839.6088
$ws.Range("N138").Value = -17119.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H52").Value = 45780
$ws.Range("J52").Value = 45780
$ws.Range("L52").Value = 45780
$ws.Range("N52").Value = -46416

$ws.Range("H61").Value = 2323.3635
$ws.Range("I61").Value = 1955.7
$ws.Range("K61").Value = 1955.7
$ws.Range("M61").Value = -1743.7

$ws.Range("H74").Value = 142858350
$ws.Range("I74").Value = 200000980
$ws.Range("J74").Value = 1757
$ws.Range("K74").Value = 200000980
$ws.Range("L74").Value = 1757
$ws.Range("M74").Value = -200000106
$ws.Range("N74").Value = -3505

$ws.Range("H77").Value = 142858350
$ws.Range("I77").Value = 200000980
$ws.Range("J77").Value = 1757
$ws.Range("K77").Value = 1000004900
$ws.Range("L77").Value = 8785
$ws.Range("M77").Value = -1000000532
$ws.Range("N77").Value = -17521

$ws.Range("H136").Value = 2323.3635
$ws.Range("I136").Value = 1955.7
$ws.Range("K136").Value = 5867.1
$ws.Range("M136").Value = -3317.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2551
$ws.Range("I20").Value = 3913.25
$ws.Range("K20").Value = 3913.25
$ws.Range("M20").Value = -3666.25

$ws.Range("H99").Value = 2085.3845
$ws.Range("I99").Value = 1739.8
$ws.Range("J99").Value = 2301.375
$ws.Range("K99").Value = 1739.8
$ws.Range("L99").Value = 2301.375
$ws.Range("M99").Value = -241.8
$ws.Range("N99").Value = -5297.375

$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

$ws.Range("H112").Value = 45500
$ws.Range("J112").Value = 45500
$ws.Range("L112").Value = 45500
$ws.Range("N112").Value = -48454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12408.317
$ws.Range("I31").Value = 19663.363
$ws.Range("K31").Value = 19663.363
$ws.Range("M31").Value = -19368.363

$ws.Range("H34").Value = 12408.317
$ws.Range("I34").Value = 19663.363
$ws.Range("K34").Value = 19663.363
$ws.Range("M34").Value = -19461.363

$ws.Range("H52").Value = 38749.5
$ws.Range("J52").Value = 38749.5
$ws.Range("L52").Value = 38749.5
$ws.Range("N52").Value = -39337.5

$ws.Range("H135").Value = 45774.8
$ws.Range("J135").Value = 45774.8
$ws.Range("L135").Value = 45774.8
$ws.Range("N135").Value = -55914.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 143660.6
$ws.Range("J131").Value = 147866.31
$ws.Range("L131").Value = 443598.93
$ws.Range("N131").Value = -453678.93

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5343.4736
$ws.Range("I70").Value = 5359.75
$ws.Range("J70").Value = 5315.5713
$ws.Range("K70").Value = 5359.75
$ws.Range("L70").Value = 5315.5713
$ws.Range("M70").Value = -5089.75
$ws.Range("N70").Value = -5855.5713

$ws.Range("H73").Value = 5343.4736
$ws.Range("I73").Value = 5359.75
$ws.Range("J73").Value = 5315.5713
$ws.Range("K73").Value = 5359.75
$ws.Range("L73").Value = 5315.5713
$ws.Range("M73").Value = -4423.75
$ws.Range("N73").Value = -7187.5713

$ws.Range("H80").Value = 3357.7083
$ws.Range("I80").Value = 3055
$ws.Range("K80").Value = 3055
$ws.Range("M80").Value = -2057

$ws.Range("H83").Value = 3357.7083
$ws.Range("I83").Value = 3055
$ws.Range("K83").Value = 15275
$ws.Range("M83").Value = -10283

$ws.Range("H132").Value = 17552.344
$ws.Range("I132").Value = 3162.3215
$ws.Range("K132").Value = 9486.9645
$ws.Range("M132").Value = -6956.9645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2690.8125
$ws.Range("I82").Value = 2522.7273
$ws.Range("J82").Value = 3060.6
$ws.Range("K82").Value = 2522.7273
$ws.Range("L82").Value = 3060.6
$ws.Range("M82").Value = -2161.7273
$ws.Range("N82").Value = -3782.6

$ws.Range("H85").Value = 2690.8125
$ws.Range("I85").Value = 2522.7273
$ws.Range("J85").Value = 3060.6
$ws.Range("K85").Value = 2522.7273
$ws.Range("L85").Value = 3060.6
$ws.Range("M85").Value = -1274.7273
$ws.Range("N85").Value = -5556.6

$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314

$ws.Range("H122").Value = 1034367.44
$ws.Range("I122").Value = 1155487.1
$ws.Range("J122").Value = 4850
$ws.Range("K122").Value = 3466461.3
$ws.Range("L122").Value = 14550
$ws.Range("M122").Value = -3464011.3
$ws.Range("N122").Value = -19450

$ws.Range("H140").Value = 48943.2
$ws.Range("J140").Value = 48943.2
$ws.Range("L140").Value = 48943.2
$ws.Range("N140").Value = -59303.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 30304898
$ws.Range("I136").Value = 43479852
$ws.Range("J136").Value = 2510.4
$ws.Range("K136").Value = 130439556
$ws.Range("L136").Value = 7531.200000000001
$ws.Range("M136").Value = -130437006
$ws.Range("N136").Value = -12631.2
